# Update optimization result values in columns A:F for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 50,6

$data[0,0] = 23.875006942982093
$data[0,1] = 109.41263880710096
$data[0,2] = 0.4697217659596237
$data[0,3] = -67439.88852573211
$data[0,4] = 1792.7814919636269
$data[0,5] = -144361.53111926137

$data[1,0] = 23.875006942982093
$data[1,1] = 109.41263880710096
$data[1,2] = 0.4697217659596237
$data[1,3] = -67439.88852573211
$data[1,4] = 1792.7814919636269
$data[1,5] = -144361.53111926137

$data[2,0] = 25.61815391394012
$data[2,1] = 119.92344667163808
$data[2,2] = 0.4912862302271642
$data[2,3] = -79189.89396393804
$data[2,4] = 1978.638666954326
$data[2,5] = -186531.45901636404

$data[3,0] = 23.874947542822227
$data[3,1] = 119.92290742064648
$data[3,2] = 0.49128622530165833
$data[3,3] = -73823.7288461798
$data[3,4] = 1938.4851287272581
$data[3,5] = -173461.91217714103

$data[4,0] = 25.224459513620715
$data[4,1] = 100.01955122973901
$data[4,2] = 0.491286390159747
$data[4,3] = -65205.51278652764
$data[4,4] = 1692.9077353969851
$data[4,5] = -127647.51180694136

$data[5,0] = 25.712531550521202
$data[5,1] = 109.41261405750602
$data[5,2] = 0.4912861309157726
$data[5,3] = -72604.9844441969
$data[5,4] = 1834.7116123739713
$data[5,5] = -155825.55759322413

$data[6,0] = 23.875006942982093
$data[6,1] = 109.41263880710096
$data[6,2] = 0.4697217659596237
$data[6,3] = -67439.88852573211
$data[6,4] = 1792.7814919636269
$data[6,5] = -144361.53111926137

$data[7,0] = 25.854488806783635
$data[7,1] = 113.52638339802725
$data[7,2] = 0.49128768736562645
$data[7,3] = -75709.77360533188
$data[7,4] = 1895.162253833972
$data[7,5] = -168730.54268168294

$data[8,0] = 25.22446318232088
$data[8,1] = 109.41251797000442
$data[8,2] = 0.49128667225964967
$data[8,3] = -71233.08613754337
$data[8,4] = 1823.4700525046828
$data[8,5] = -152780.27504011022

$data[9,0] = 25.60112650573445
$data[9,1] = 100.01981785075569
$data[9,2] = 0.49128640717886446
$data[9,3] = -66174.3938530885
$data[9,4] = 1701.5860819514926
$data[9,5] = -129611.64249073368

$data[10,0] = 23.874944903946343
$data[10,1] = 119.92337605130746
$data[10,2] = 0.4912876684190492
$data[10,3] = -73824.0053793126
$data[10,4] = 1938.4915559728834
$data[10,5] = -173463.2494630023

$data[11,0] = 23.874947542822227
$data[11,1] = 119.92290742064648
$data[11,2] = 0.49128622530165833
$data[11,3] = -73823.7288461798
$data[11,4] = 1938.4851287272581
$data[11,5] = -173461.91217714103

$data[12,0] = 25.224455268661256
$data[12,1] = 124.89968618097143
$data[12,2] = 0.49128623968856905
$data[12,3] = -81171.33680073943
$data[12,4] = 2038.7415161631711
$data[12,5] = -199147.5218913121

$data[13,0] = 25.224461090432516
$data[13,1] = 132.16657154812918
$data[13,2] = 0.4912861045693819
$data[13,3] = -85834.59057352405
$data[13,4] = 2139.751359271499
$data[13,5] = -223018.88539772257

$data[14,0] = 25.618153097884537
$data[14,1] = 124.89982348896177
$data[14,2] = 0.4912875358984885
$data[14,3] = -82433.12470589789
$data[14,4] = 2047.8102624453945
$data[14,5] = -202349.5850209436

$data[15,0] = 25.61815391394012
$data[15,1] = 119.92344667163808
$data[15,2] = 0.4912862302271642
$data[15,3] = -79189.89396393804
$data[15,4] = 1978.638666954326
$data[15,5] = -186531.45901636404

$data[16,0] = 25.601125681965666
$data[16,1] = 100.01936674002421
$data[16,2] = 0.49128624323855974
$data[16,3] = -66174.09792718894
$data[16,4] = 1701.5797954885766
$data[16,5] = -129610.46762477436

$data[17,0] = 25.85448110476594
$data[17,1] = 132.166496608159
$data[17,2] = 0.49128621147357765
$data[17,3] = -87970.06190036972
$data[17,4] = 2154.259676613875
$data[17,5] = -228756.1409334963

$data[18,0] = 25.854488806783635
$data[18,1] = 113.52638339802725
$data[18,2] = 0.49128768736562645
$data[18,3] = -75709.77360533188
$data[18,4] = 1895.162253833972
$data[18,5] = -168730.54268168294

$data[19,0] = 25.712525614235773
$data[19,1] = 113.52603542201737
$data[19,2] = 0.4912873592801911
$data[19,3] = -75295.6663194635
$data[19,4] = 1891.8880105420335
$data[19,5] = -167775.89293842294

$data[20,0] = 25.61814424404588
$data[20,1] = 113.52623549877508
$data[20,2] = 0.49128752114779695
$data[20,3] = -75020.63799301833
$data[20,4] = 1889.717185743113
$data[20,5] = -167142.4914544365

$data[21,0] = 25.224461090432516
$data[21,1] = 132.16657154812918
$data[21,2] = 0.4912861045693819
$data[21,3] = -85834.59057352405
$data[21,4] = 2139.751359271499
$data[21,5] = -223018.88539772257

$data[22,0] = 25.224451821106776
$data[22,1] = 132.1658771138473
$data[22,2] = 0.49128764059098046
$data[22,3] = -85834.1135496532
$data[22,4] = 2139.7414655447405
$data[22,5] = -223016.45526407374

$data[23,0] = 23.875006942982093
$data[23,1] = 109.41263880710096
$data[23,2] = 0.4697217659596237
$data[23,3] = -67439.88852573211
$data[23,4] = 1792.7814919636269
$data[23,5] = -144361.53111926137

$data[24,0] = 25.60112650573445
$data[24,1] = 100.01981785075569
$data[24,2] = 0.49128640717886446
$data[24,3] = -66174.3938530885
$data[24,4] = 1701.5860819514926
$data[24,5] = -129611.64249073368

$data[25,0] = 25.601125681965666
$data[25,1] = 100.01936674002421
$data[25,2] = 0.49128624323855974
$data[25,3] = -66174.09792718894
$data[25,4] = 1701.5797954885766
$data[25,5] = -129610.46762477436

$data[26,0] = 23.874944903946343
$data[26,1] = 119.92337605130746
$data[26,2] = 0.4912876684190492
$data[26,3] = -73824.0053793126
$data[26,4] = 1938.4915559728834
$data[26,5] = -173463.2494630023

$data[27,0] = 25.85449174397755
$data[27,1] = 109.41223618051394
$data[27,2] = 0.49128628533875596
$data[27,3] = -73003.75007605697
$data[27,4] = 1837.975700362556
$data[27,5] = -156710.1505955151

$data[28,0] = 25.618153097884537
$data[28,1] = 124.89982348896177
$data[28,2] = 0.4912875358984885
$data[28,3] = -82433.12470589789
$data[28,4] = 2047.8102624453945
$data[28,5] = -202349.5850209436

$data[29,0] = 25.618149936820124
$data[29,1] = 119.92340383371194
$data[29,2] = 0.49128638615123266
$data[29,3] = -79189.85380534551
$data[29,4] = 1978.6379771105644
$data[29,5] = -186531.2958052628

$data[30,0] = 25.85445982289632
$data[30,1] = 119.92350664750568
$data[30,2] = 0.4912876814033579
$data[30,3] = -79917.31668954044
$data[30,4] = 1984.0815996099989
$data[30,5] = -188303.12004191964

$data[31,0] = 25.618152863259823
$data[31,1] = 136.1051215133184
$data[31,2] = 0.4912861457984144
$data[31,3] = -89735.90585324344
$data[31,4] = 2203.563924574544
$data[31,5] = -240323.62677455996

$data[32,0] = 25.61815323872265
$data[32,1] = 136.10530457568424
$data[32,2] = 0.49128626646661977
$data[32,3] = -89736.02647165253
$data[32,4] = 2203.566475618724
$data[32,5] = -240324.2774453686

$data[33,0] = 25.854490710176915
$data[33,1] = 136.10522450456983
$data[33,2] = 0.49128762983616314
$data[33,3] = -90560.74644715297
$data[33,4] = 2209.0081900844407
$data[33,5] = -242606.59284989868

$data[34,0] = 25.712531849715972
$data[34,1] = 136.10500941319125
$data[34,2] = 0.4912861536921094
$data[34,3] = -90065.19671759568
$data[34,4] = 2205.735914298933
$data[34,5] = -241234.76033456274

$data[35,0] = 25.618141575828098
$data[35,1] = 136.1052892094993
$data[35,2] = 0.49128626220529353
$data[35,3] = -89735.9757557857
$data[35,4] = 2203.56599350891
$data[35,5] = -240324.11048978425

$data[36,0] = 25.22444665366699
$data[36,1] = 134.21966882959916
$data[36,2] = 0.49128607181745276
$data[36,3] = -87152.03453624238
$data[36,4] = 2168.289079594101
$data[36,5] = -230007.79738395623

$data[37,0] = 25.854469424339616
$data[37,1] = 132.16648355471114
$data[37,2] = 0.4912874839933066
$data[37,3] = -87970.01373995544
$data[37,4] = 2154.2592032908265
$data[37,5] = -228755.989330498

$data[38,0] = 25.854476872103522
$data[38,1] = 132.16634503828544
$data[38,2] = 0.4912860688119088
$data[38,3] = -87969.94785811675
$data[38,4] = 2154.2574748794736
$data[38,5] = -228755.577231372

$data[39,0] = 25.22446484632068
$data[39,1] = 132.16576405889214
$data[39,2] = 0.4912876827691401
$data[39,3] = -85834.08515182372
$data[39,4] = 2139.740193293177
$data[39,5] = -223016.1919967712

$data[40,0] = 25.618154901458492
$data[40,1] = 124.89955852928965
$data[40,2] = 0.4912865957801983
$data[40,3] = -82432.95779198324
$data[40,4] = 2047.8066379455877
$data[40,5] = -202348.7403435061

$data[41,0] = 25.224457167100017
$data[41,1] = 124.89979492531337
$data[41,2] = 0.4912874625406506
$data[41,3] = -81171.41268378262
$data[41,4] = 2038.7430494436885
$data[41,5] = -199147.88443845845

$data[42,0] = 25.712526358825286
$data[42,1] = 119.9234638988224
$data[42,2] = 0.491287475222759
$data[42,3] = -79480.39708154938
$data[42,4] = 1980.8122814328726
$data[42,5] = -187238.9779896054

$data[43,0] = 25.224454393768845
$data[43,1] = 119.9229672408103
$data[43,2] = 0.49128748426137475
$data[43,3] = -77977.72310673723
$data[43,4] = 1969.5650803687404
$data[43,5] = -183578.61250759405

$data[44,0] = 25.22447538390697
$data[44,1] = 136.10557606294233
$data[44,2] = 0.491287408411703
$data[44,3] = -88362.34020298117
$data[44,4] = 2194.5038277630333
$data[44,5] = -236522.99149003005

$data[45,0] = 25.85447205251505
$data[45,1] = 132.1665110141271
$data[45,2] = 0.49128736098724024
$data[45,3] = -87970.04070784891
$data[45,4] = 2154.2596477152383
$data[45,5] = -228756.1084061491

$data[46,0] = 25.22449257972883
$data[46,1] = 132.16640075832268
$data[46,2] = 0.49128638480306946
$data[46,3] = -85834.58771599362
$data[46,4] = 2139.7497054530813
$data[46,5] = -223018.59525675667

$data[47,0] = 24.85273026480196
$data[47,1] = 132.16639478769378
$data[47,2] = 0.49128762426799916
$data[47,3] = -84574.45721252587
$data[47,4] = 2131.1879140629944
$data[47,5] = -219632.99023224477

$data[48,0] = 25.854493352718652
$data[48,1] = 124.89977064302131
$data[48,2] = 0.4912874721890247
$data[48,3] = -83190.4993712753
$data[48,4] = 2053.2524451011486
$data[48,5] = -204271.3754316238

$data[49,0] = 25.712527513888848
$data[49,1] = 124.8997707999258
$data[49,2] = 0.4912874560972334
$data[49,3] = -82735.53584791279
$data[49,4] = 2049.9829743032005
$data[49,5] = -203116.8839992683

$ws.Range("A2:F51").Value = $data
